# Weekly update: insert a new daily price record as row 12 on the
# "Macroferia Regional de Talca - Poroto verde" sheet, shifting every
# existing record (old rows 12-120) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 12; this pushes the old
# row 12 (and everything below it) down to row 13, old row 120 ends up
# at row 121, and the sheet's dimension grows to A1:R121 automatically.
$ws.Rows.Item(12).Insert()

# Populate the newly-inserted row 12 with the new weekly record.
$ws.Cells.Item(12, 1).Value = 5
$ws.Cells.Item(12, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(12, 3).Value = "Maule"
$ws.Cells.Item(12, 4).Value = 44545
$ws.Cells.Item(12, 5).Value = 7
$ws.Cells.Item(12, 6).Value = 100112031
$ws.Cells.Item(12, 7).Value = "Poroto verde"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 300
$ws.Cells.Item(12, 11).Value = 12000
$ws.Cells.Item(12, 12).Value = 12000
$ws.Cells.Item(12, 13).Value = 12000
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Región del Maule"
$ws.Cells.Item(12, 16).Value = 480
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
